$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.83%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '38.73'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.38%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.093'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.92%'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.14%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.971'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.34%'

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.202'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.14%'

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '7.935'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.98%'

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9295'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.88%'

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1439'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '12.36%'

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1960'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.31%'

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09089'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.27%'

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03510'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.78%'

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09811'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.44%'

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001408'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.19%'

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005926'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-5.16%'

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.601'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-5.15%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3447'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.03%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1335'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.00%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.823'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-6.65%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2403'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-7.74%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04436'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.07%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.56%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004848'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '5.05%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001302'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '3.98%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02093'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '7.08%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05111'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-6.31%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007464'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.41%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01014'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.14%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1363'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002143'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.49%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01047'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '6.47%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006237'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.70%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.02%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003068'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-3.52%'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.02%'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.02%'
